# Rename the first sheet from "Sheet1" to "data"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "data"

# Move the selection on that sheet from A4 to B3
$ws1.Range("B3").Select()
